$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix F93: CodigoNiv1 should be text "04" (shared string) instead of numeric 4 ---
$ws.Range("F93").Value = "04"

# --- Copy row 161 formatting down through new rows 162:173 (columns B:S and V:W only, to avoid introducing stray T/U cells) ---
$ws.Range("B161:S161").Copy() | Out-Null
$ws.Range("B162:S173").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("V161:W161").Copy() | Out-Null
$ws.Range("V162:W173").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 161
$ws.Range("B161").Value = "HND"
$ws.Range("C161").Value = "Honduras"
$ws.Range("D161").Value = 3
$ws.Range("E161").Value = 2
$ws.Range("F161").Value = "02"
$ws.Range("G161").Value = "Colón"
$ws.Range("H161").Value = "Departamento"
$ws.Range("I161").Value = 3
$ws.Range("J161").Value = "0203"
$ws.Range("K161").Value = "Iriona"
$ws.Range("L161").Value = "Municipio"
$ws.Range("M161").Value = 9
$ws.Range("N161").Value = "020309"
$ws.Range("O161").Value = "Sangrelaya"
$ws.Range("P161").Value = "Aldea"
$ws.Range("Q161").Value = "HND-0203"
$ws.Range("R161").Value = "Salud"
$ws.Range("S161").Value = "Cesamo"
$ws.Range("V161").Value = 15.966279
$ws.Range("W161").Value = -85.09611

# Row 162
$ws.Range("B162").Value = "HND"
$ws.Range("C162").Value = "Honduras"
$ws.Range("D162").Value = 3
$ws.Range("E162").Value = 2
$ws.Range("F162").Value = "02"
$ws.Range("G162").Value = "Colón"
$ws.Range("H162").Value = "Departamento"
$ws.Range("I162").Value = 4
$ws.Range("J162").Value = "0204"
$ws.Range("K162").Value = "Limón"
$ws.Range("L162").Value = "Municipio"
$ws.Range("M162").Value = 1
$ws.Range("N162").Value = "020401"
$ws.Range("O162").Value = "Limón"
$ws.Range("P162").Value = "Aldea"
$ws.Range("Q162").Value = "HND-0204"
$ws.Range("R162").Value = "Salud"
$ws.Range("S162").Value = "Cesamo"
$ws.Range("V162").Value = 15.86447408
$ws.Range("W162").Value = -85.506745

# Row 163
$ws.Range("B163").Value = "HND"
$ws.Range("C163").Value = "Honduras"
$ws.Range("D163").Value = 3
$ws.Range("E163").Value = 2
$ws.Range("F163").Value = "02"
$ws.Range("G163").Value = "Colón"
$ws.Range("H163").Value = "Departamento"
$ws.Range("I163").Value = 10
$ws.Range("J163").Value = "0210"
$ws.Range("K163").Value = "Bonito Oriental"
$ws.Range("L163").Value = "Municipio"
$ws.Range("M163").Value = 1
$ws.Range("N163").Value = "021001"
$ws.Range("O163").Value = "Bonito Oriental"
$ws.Range("P163").Value = "Aldea"
$ws.Range("Q163").Value = "HND-0210"
$ws.Range("R163").Value = "Salud"
$ws.Range("S163").Value = "Cesamo"
$ws.Range("V163").Value = 15.748343
$ws.Range("W163").Value = -85.73515

# Row 164
$ws.Range("B164").Value = "HND"
$ws.Range("C164").Value = "Honduras"
$ws.Range("D164").Value = 3
$ws.Range("E164").Value = 2
$ws.Range("F164").Value = "02"
$ws.Range("G164").Value = "Colón"
$ws.Range("H164").Value = "Departamento"
$ws.Range("I164").Value = 5
$ws.Range("J164").Value = "0205"
$ws.Range("K164").Value = "Sabá"
$ws.Range("L164").Value = "Municipio"
$ws.Range("M164").Value = 11
$ws.Range("N164").Value = "020511"
$ws.Range("O164").Value = "Elixir"
$ws.Range("P164").Value = "Aldea"
$ws.Range("Q164").Value = "HND-0205"
$ws.Range("R164").Value = "Salud"
$ws.Range("S164").Value = "Cesamo"
$ws.Range("V164").Value = 15.534311
$ws.Range("W164").Value = -86.273985

# Row 165
$ws.Range("B165").Value = "HND"
$ws.Range("C165").Value = "Honduras"
$ws.Range("D165").Value = 3
$ws.Range("E165").Value = 3
$ws.Range("F165").Value = "03"
$ws.Range("G165").Value = "Comayagua"
$ws.Range("H165").Value = "Departamento"
$ws.Range("I165").Value = 20
$ws.Range("J165").Value = "0320"
$ws.Range("K165").Value = "Las Lajas"
$ws.Range("L165").Value = "Municipio"
$ws.Range("M165").Value = 1
$ws.Range("N165").Value = "032001"
$ws.Range("O165").Value = "Las Lajas"
$ws.Range("P165").Value = "Aldea"
$ws.Range("Q165").Value = "HND-0320"
$ws.Range("R165").Value = "Salud"
$ws.Range("S165").Value = "Cesamo"
$ws.Range("V165").Value = 14.894002
$ws.Range("W165").Value = -87.579126

# Row 166
$ws.Range("B166").Value = "HND"
$ws.Range("C166").Value = "Honduras"
$ws.Range("D166").Value = 3
$ws.Range("E166").Value = 3
$ws.Range("F166").Value = "03"
$ws.Range("G166").Value = "Comayagua"
$ws.Range("H166").Value = "Departamento"
$ws.Range("I166").Value = 21
$ws.Range("J166").Value = "0321"
$ws.Range("K166").Value = "Taulabé"
$ws.Range("L166").Value = "Municipio"
$ws.Range("M166").Value = 1
$ws.Range("N166").Value = "032101"
$ws.Range("O166").Value = "Taulabé"
$ws.Range("P166").Value = "Aldea"
$ws.Range("Q166").Value = "HND-0321"
$ws.Range("R166").Value = "Salud"
$ws.Range("S166").Value = "Cesamo"
$ws.Range("V166").Value = 14.690765
$ws.Range("W166").Value = -87.965248

# Row 167
$ws.Range("B167").Value = "HND"
$ws.Range("C167").Value = "Honduras"
$ws.Range("D167").Value = 3
$ws.Range("E167").Value = 3
$ws.Range("F167").Value = "03"
$ws.Range("G167").Value = "Comayagua"
$ws.Range("H167").Value = "Departamento"
$ws.Range("I167").Value = 3
$ws.Range("J167").Value = "0303"
$ws.Range("K167").Value = "El Rosario"
$ws.Range("L167").Value = "Municipio"
$ws.Range("M167").Value = 1
$ws.Range("N167").Value = "030301"
$ws.Range("O167").Value = "El Rosario"
$ws.Range("P167").Value = "Aldea"
$ws.Range("Q167").Value = "HND-0303"
$ws.Range("R167").Value = "Salud"
$ws.Range("S167").Value = "Cesamo"
$ws.Range("V167").Value = 14.575887
$ws.Range("W167").Value = -87.72896

# Row 168
$ws.Range("B168").Value = "HND"
$ws.Range("C168").Value = "Honduras"
$ws.Range("D168").Value = 3
$ws.Range("E168").Value = 3
$ws.Range("F168").Value = "03"
$ws.Range("G168").Value = "Comayagua"
$ws.Range("H168").Value = "Departamento"
$ws.Range("I168").Value = 19
$ws.Range("J168").Value = "0319"
$ws.Range("K168").Value = "Villa de San Antonio"
$ws.Range("L168").Value = "Municipio"
$ws.Range("M168").Value = 1
$ws.Range("N168").Value = "031901"
$ws.Range("O168").Value = "Villa de San Antonio"
$ws.Range("P168").Value = "Aldea"
$ws.Range("Q168").Value = "HND-0319"
$ws.Range("R168").Value = "Salud"
$ws.Range("S168").Value = "Cesamo"
$ws.Range("V168").Value = 14.323913
$ws.Range("W168").Value = -87.61369

# Row 169
$ws.Range("B169").Value = "HND"
$ws.Range("C169").Value = "Honduras"
$ws.Range("D169").Value = 3
$ws.Range("E169").Value = 4
$ws.Range("F169").Value = "04"
$ws.Range("G169").Value = "Copán"
$ws.Range("H169").Value = "Departamento"
$ws.Range("I169").Value = 6
$ws.Range("J169").Value = "06"
$ws.Range("K169").Value = "Cucuyagua"
$ws.Range("L169").Value = "Municipio"
$ws.Range("O169").Value = "Cucuyagua"
$ws.Range("P169").Value = "Aldea"
$ws.Range("Q169").Value = "HND-0406"
$ws.Range("R169").Value = "Salud"
$ws.Range("S169").Value = "Cesamo"
$ws.Range("V169").Value = 14.645818
$ws.Range("W169").Value = -88.873899

# Row 170
$ws.Range("B170").Value = "HND"
$ws.Range("C170").Value = "Honduras"
$ws.Range("D170").Value = 3
$ws.Range("E170").Value = 4
$ws.Range("F170").Value = "04"
$ws.Range("G170").Value = "Copán"
$ws.Range("H170").Value = "Departamento"
$ws.Range("I170").Value = 8
$ws.Range("J170").Value = "08"
$ws.Range("K170").Value = "Dulce Nombre"
$ws.Range("L170").Value = "Municipio"
$ws.Range("O170").Value = "Dulce Nombre"
$ws.Range("P170").Value = "Aldea"
$ws.Range("Q170").Value = "HND-0408"
$ws.Range("R170").Value = "Salud"
$ws.Range("S170").Value = "Cesamo"
$ws.Range("V170").Value = 14.845943
$ws.Range("W170").Value = -88.83152

# Row 171
$ws.Range("B171").Value = "HND"
$ws.Range("C171").Value = "Honduras"
$ws.Range("D171").Value = 3
$ws.Range("E171").Value = 4
$ws.Range("F171").Value = "04"
$ws.Range("G171").Value = "Copán"
$ws.Range("H171").Value = "Departamento"
$ws.Range("I171").Value = 19
$ws.Range("J171").Value = "19"
$ws.Range("K171").Value = "San Nicolas"
$ws.Range("L171").Value = "Municipio"
$ws.Range("O171").Value = "San Nicolas"
$ws.Range("P171").Value = "Aldea"
$ws.Range("Q171").Value = "HND-0419"
$ws.Range("R171").Value = "Salud"
$ws.Range("S171").Value = "Cesamo"
$ws.Range("V171").Value = 15.001558
$ws.Range("W171").Value = -88.751771

# Row 172
$ws.Range("B172").Value = "HND"
$ws.Range("C172").Value = "Honduras"
$ws.Range("D172").Value = 3
$ws.Range("E172").Value = 4
$ws.Range("F172").Value = "04"
$ws.Range("G172").Value = "Copán"
$ws.Range("H172").Value = "Departamento"
$ws.Range("I172").Value = 13
$ws.Range("J172").Value = "13"
$ws.Range("K172").Value = "Nueva Arcadia"
$ws.Range("L172").Value = "Municipio"
$ws.Range("M172").Value = 1
$ws.Range("N172").Value = "041301"
$ws.Range("O172").Value = "La Entrada"
$ws.Range("P172").Value = "Aldea"
$ws.Range("Q172").Value = "HND-0413"
$ws.Range("R172").Value = "Salud"
$ws.Range("S172").Value = "Cesamo"
$ws.Range("V172").Value = 15.061814
$ws.Range("W172").Value = -88.746099

# Row 173
$ws.Range("B173").Value = "HND"
$ws.Range("C173").Value = "Honduras"
$ws.Range("D173").Value = 3
$ws.Range("E173").Value = 4
$ws.Range("F173").Value = "04"
$ws.Range("G173").Value = "Copán"
$ws.Range("H173").Value = "Departamento"
$ws.Range("I173").Value = 10
$ws.Range("J173").Value = "10"
$ws.Range("K173").Value = "Florida"
$ws.Range("L173").Value = "Municipio"
$ws.Range("M173").Value = 1
$ws.Range("N173").Value = "041001"
$ws.Range("O173").Value = "Florida"
$ws.Range("P173").Value = "Aldea"
$ws.Range("Q173").Value = "HND-0410"
$ws.Range("R173").Value = "Salud"
$ws.Range("S173").Value = "Cesamo"
$ws.Range("V173").Value = 15.024939
$ws.Range("W173").Value = -88.835689

# --- Resize table (ListObject) to include new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:W173"))

# --- Update _xlnm._FilterDatabase defined name range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "HOSPITALES!_FilterDatabase") {
        $n.RefersTo = "=HOSPITALES!`$A`$1:`$W`$173"
    }
}

# --- Update view: scroll/pane position and active selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 162
$win.ScrollColumn = 21
$ws.Range("W174").Select()
